$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(-18.6469453757543, 1.948105357729458, -18.6469453757543, -18.6469453757543, -18.6469453757543, -18.6469453757543, -18.6469453757543, -18.6469453757543, -18.6469453757543, -18.6469453757543),
    @(-18.6469453757543, -18.6469453757543, -18.6469453757543, -18.6469453757543, -18.6469453757543, -18.6469453757543, -18.6469453757543, 1.218801485533855, -18.6469453757543, -18.6469453757543),
    @(-18.6469453757543, 1.999093809060723, 1.652308486473567, -18.6469453757543, 3.424112662232885, -18.6469453757543, -18.6469453757543, -18.6469453757543, 1.000043932716085, -18.6469453757543),
    @(-18.6469453757543, 1.653493212068003, -18.6469453757543, -18.6469453757543, -18.6469453757543, 2.750650403607153, -18.6469453757543, -18.6469453757543, -18.6469453757543, -18.6469453757543),
    @(-18.6469453757543, -18.6469453757543, -18.6469453757543, -18.6469453757543, -18.6469453757543, -18.6469453757543, -18.6469453757543, -18.6469453757543, -18.6469453757543, -18.6469453757543),
    @(2.440955732114179, -18.6469453757543, -18.6469453757543, -18.6469453757543, -18.6469453757543, -18.6469453757543, -18.6469453757543, -18.6469453757543, -18.6469453757543, -18.6469453757543),
    @(-18.6469453757543, -18.6469453757543, -18.6469453757543, 1.774084949643362, -18.6469453757543, -18.6469453757543, -18.6469453757543, -18.6469453757543, -18.6469453757543, -18.6469453757543),
    @(3.864923313103878, -18.6469453757543, -18.6469453757543, -18.6469453757543, -18.6469453757543, -18.6469453757543, -18.6469453757543, -18.6469453757543, -18.6469453757543, -18.6469453757543),
    @(-18.6469453757543, -18.6469453757543, -18.6469453757543, -18.6469453757543, -18.6469453757543, -18.6469453757543, -18.6469453757543, 1.521380294553261, -18.6469453757543, 2.189620500958363),
    @(-18.6469453757543, -18.6469453757543, -18.6469453757543, 2.916612616220791, -18.6469453757543, 2.888770921104235, -18.6469453757543, -18.6469453757543, -18.6469453757543, 1.877834480920477),
    @(-18.6469453757543, -18.6469453757543, -18.6469453757543, -18.6469453757543, -18.6469453757543, -18.6469453757543, -18.6469453757543, -18.6469453757543, -18.6469453757543, -18.6469453757543),
    @(-18.6469453757543, -18.6469453757543, -18.6469453757543, 2.530995403022448, -18.6469453757543, -18.6469453757543, -18.6469453757543, -18.6469453757543, 1.710126541219667, 1.832350380455183),
    @(-18.6469453757543, -18.6469453757543, 1.571387265452801, -18.6469453757543, -18.6469453757543, -18.6469453757543, -18.6469453757543, -18.6469453757543, -18.6469453757543, 1.938949657413982),
    @(-18.6469453757543, -18.6469453757543, 1.740652338236284, -18.6469453757543, -18.6469453757543, -18.6469453757543, -18.6469453757543, -18.6469453757543, -18.6469453757543, -18.6469453757543),
    @(-18.6469453757543, -18.6469453757543, -18.6469453757543, -18.6469453757543, -18.6469453757543, -18.6469453757543, -18.6469453757543, -18.6469453757543, 1.995960386309749, -18.6469453757543),
    @(-18.6469453757543, 2.151476587564451, 1.878102223768725, -18.6469453757543, -18.6469453757543, -18.6469453757543, -18.6469453757543, 1.950107835708667, 2.525475418100582, -18.6469453757543),
    @(-18.6469453757543, -18.6469453757543, -18.6469453757543, -18.6469453757543, -18.6469453757543, -18.6469453757543, -18.6469453757543, 2.174922788062998, 2.316597751301674, -18.6469453757543),
    @(-18.6469453757543, -18.6469453757543, 2.030885217518834, -18.6469453757543, -18.6469453757543, -18.6469453757543, -18.6469453757543, 2.022380033216714, -18.6469453757543, -18.6469453757543),
    @(-18.6469453757543, 1.077130425131631, 1.475277241458547, -18.6469453757543, 3.211943591763764, -18.6469453757543, -18.6469453757543, 1.238538011348082, -18.6469453757543, 2.126790268499105),
    @(-18.6469453757543, 1.289952751641185, -18.6469453757543, 1.700144996410783, -18.6469453757543, 2.551728066338604, 4.321924755946238, -18.6469453757543, -18.6469453757543, -18.6469453757543)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $rowvals = $data[$i]
    for ($j = 0; $j -lt $rowvals.Length; $j++) {
        $col = $j + 2
        $ws.Cells.Item($row, $col).Value = $rowvals[$j]
    }
}
